$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set the uniform fill value across the full data range B2:K21
$ws.Range("B2:K21").Value = -18.35557047814295

# Step 2: apply the specific (non-uniform) values that differ from the fill
$cellValues = @{
    "C2" = 2.474199476936509
    "I3" = 2.379870682204503
    "C4" = 2.166797454256023
    "D4" = 2.898757515923904
    "F4" = 2.492169731024573
    "H4" = 1.70898416032878
    "C5" = 0.9347566694192283
    "G5" = 2.244481248420171
    "B7" = 2.978458878643612
    "E8" = 2.890022265798286
    "B9" = 3.599140223856239
    "I10" = 1.374489955405986
    "K10" = 2.160711465050086
    "E11" = 1.92737878914738
    "G11" = 2.350417038721486
    "K11" = 1.309857395150263
    "E13" = 1.558524088491349
    "K13" = 1.27540935615156
    "D14" = 1.657611777386651
    "K14" = 1.995430192343434
    "D15" = -0.3415663585807351
    "C17" = 0.5671445723064817
    "D17" = -0.1802304781181286
    "H17" = 0.2435759902290854
    "I17" = 0.7538036516848216
    "J17" = 4.321924008670065
    "H18" = -0.07182645432115185
    "I18" = 0.3816405487847508
    "D19" = 1.651491017739423
    "H19" = 1.747667378137802
    "I19" = 2.264546407595732
    "C20" = 1.671451322598567
    "D20" = 2.193410367513295
    "F20" = 3.845360745074863
    "H20" = 2.510945544154019
    "I20" = 2.140191856169461
    "K20" = 2.731488936574111
    "C21" = 1.754720379974387
    "E21" = 2.545468008679666
    "G21" = 3.345050842247196
    "H21" = 2.469159126134352
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}

Write-Host "Applied PSSM value updates"
